$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.920312643051147
$ws.Range("B1").Value = 3.730503082275391
$ws.Range("C1").Value = 1.947177886962891
$ws.Range("D1").Value = 1.457335948944092
$ws.Range("E1").Value = 1.29633641242981
